$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 -> "Experimental" property: set its value to the literal text "false"
# (using a formula that evaluates to the string, then converting it to a
# plain value via copy/paste-special, avoids Excel's automatic TRUE/FALSE
# boolean auto-conversion while keeping the original cell style/format)
$ws.Cells.Item(7, 2).Formula = '="false"'
$ws.Cells.Item(7, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)

# Row 14 -> "Case Sensitive" property: set its value to the literal text "true"
$ws.Cells.Item(14, 2).Formula = '="true"'
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)

$excel.CutCopyMode = $false
